# Auto-generated edit script applying scheduled-runner market-price updates
# to the Jenova_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 233.875
$ws.Range("J6").Value = 207.5
$ws.Range("L6").Value = 622.5
$ws.Range("N6").Value = -846.5
$ws.Range("H96").Value = 287.6
$ws.Range("I96").Value = 252.4
$ws.Range("J96").Value = 322.8
$ws.Range("K96").Value = 757.2
$ws.Range("L96").Value = 968.4000000000001
$ws.Range("M96").Value = 615.8
$ws.Range("N96").Value = -3714.4
$ws.Range("H103").Value = 358.2857
$ws.Range("J103").Value = 600
$ws.Range("L103").Value = 1800
$ws.Range("N103").Value = -2972
$ws.Range("H132").Value = 4948.0557
$ws.Range("I132").Value = 5252.7744
$ws.Range("J132").Value = 3058.8
$ws.Range("K132").Value = 15758.3232
$ws.Range("L132").Value = 9176.400000000001
$ws.Range("M132").Value = -13228.3232
$ws.Range("N132").Value = -14236.4
$ws.Range("H137").Value = 2887.7917
$ws.Range("I137").Value = 3006.2
$ws.Range("J137").Value = 2856.6316
$ws.Range("K137").Value = 9018.599999999999
$ws.Range("L137").Value = 8569.8948
$ws.Range("M137").Value = -6468.599999999999
$ws.Range("N137").Value = -13669.8948
$ws.Range("H138").Value = 6446.907
$ws.Range("J138").Value = 7557.2104
$ws.Range("L138").Value = 22671.6312
$ws.Range("N138").Value = -32951.6312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1773.2858
$ws.Range("I32").Value = 1597.6857
$ws.Range("K32").Value = 1597.6857
$ws.Range("M32").Value = -1310.6857
$ws.Range("H61").Value = 4109.893
$ws.Range("I61").Value = 3394.6956
$ws.Range("K61").Value = 3394.6956
$ws.Range("M61").Value = -3182.6956
$ws.Range("H63").Value = 7457.4287
$ws.Range("I63").Value = 4140.8
$ws.Range("J63").Value = 9300
$ws.Range("K63").Value = 4140.8
$ws.Range("L63").Value = 9300
$ws.Range("M63").Value = -3454.8
$ws.Range("N63").Value = -10672
$ws.Range("H66").Value = 7457.4287
$ws.Range("I66").Value = 4140.8
$ws.Range("J66").Value = 9300
$ws.Range("K66").Value = 20704
$ws.Range("L66").Value = 46500
$ws.Range("M66").Value = -17272
$ws.Range("N66").Value = -53364
$ws.Range("H74").Value = 402518.2
$ws.Range("I74").Value = 626934.75
$ws.Range("K74").Value = 626934.75
$ws.Range("M74").Value = -626060.75
$ws.Range("H77").Value = 402518.2
$ws.Range("I77").Value = 626934.75
$ws.Range("K77").Value = 3134673.75
$ws.Range("M77").Value = -3130305.75
$ws.Range("H132").Value = 393077.47
$ws.Range("I132").Value = 1078827.2
$ws.Range("J132").Value = 12105.389
$ws.Range("K132").Value = 3236481.6
$ws.Range("L132").Value = 36316.167
$ws.Range("M132").Value = -3233951.6
$ws.Range("N132").Value = -41376.167
$ws.Range("H136").Value = 4109.893
$ws.Range("I136").Value = 3394.6956
$ws.Range("K136").Value = 10184.0868
$ws.Range("M136").Value = -7634.086800000001
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2496.4
$ws.Range("J20").Value = 2135.4
$ws.Range("L20").Value = 2135.4
$ws.Range("N20").Value = -2629.4
$ws.Range("H105").Value = 73895.78999999999
$ws.Range("I105").Value = 126695.75
$ws.Range("J105").Value = 3495.8333
$ws.Range("K105").Value = 126695.75
$ws.Range("L105").Value = 3495.8333
$ws.Range("M105").Value = -124948.75
$ws.Range("N105").Value = -6989.8333
$ws.Range("H134").Value = 43014
$ws.Range("I134").Value = 2254.762
$ws.Range("K134").Value = 6764.286
$ws.Range("M134").Value = -4229.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 501994.34
$ws.Range("J31").Value = 4418.84
$ws.Range("L31").Value = 4418.84
$ws.Range("N31").Value = -5008.84
$ws.Range("H34").Value = 501994.34
$ws.Range("J34").Value = 4418.84
$ws.Range("L34").Value = 4418.84
$ws.Range("N34").Value = -4822.84
$ws.Range("H122").Value = 12332
$ws.Range("I122").Value = 18040.4
$ws.Range("K122").Value = 54121.2
$ws.Range("M122").Value = -51671.2
$ws.Range("H132").Value = 2745.1428
$ws.Range("I132").Value = 2431.9583
$ws.Range("J132").Value = 4624.25
$ws.Range("K132").Value = 7295.874899999999
$ws.Range("L132").Value = 13872.75
$ws.Range("M132").Value = -4765.874899999999
$ws.Range("N132").Value = -18932.75
$ws.Range("H134").Value = 273234.16
$ws.Range("I134").Value = 2746.238
$ws.Range("J134").Value = 628249.5600000001
$ws.Range("K134").Value = 8238.714
$ws.Range("L134").Value = 1884748.68
$ws.Range("M134").Value = -5703.714
$ws.Range("N134").Value = -1889818.68

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 586464.0600000001
$ws.Range("I4").Value = 464456.75
$ws.Range("K4").Value = 1393370.25
$ws.Range("M4").Value = -1393258.25
$ws.Range("H17").Value = 6006.5
$ws.Range("I17").Value = 6006.5
$ws.Range("K17").Value = 18019.5
$ws.Range("M17").Value = -17850.5
$ws.Range("H44").Value = 171.5
$ws.Range("I44").Value = 171.5
$ws.Range("K44").Value = 514.5
$ws.Range("M44").Value = -116.5
$ws.Range("H107").Value = 25704.045
$ws.Range("I107").Value = 633.0909
$ws.Range("J107").Value = 34061.03
$ws.Range("K107").Value = 1899.2727
$ws.Range("L107").Value = 102183.09
$ws.Range("M107").Value = 20.72730000000001
$ws.Range("N107").Value = -106023.09
$ws.Range("H129").Value = 8710.77
$ws.Range("I129").Value = 1090.5
$ws.Range("J129").Value = 20903.2
$ws.Range("K129").Value = 3271.5
$ws.Range("L129").Value = 62709.60000000001
$ws.Range("M129").Value = 1728.5
$ws.Range("N129").Value = -72709.60000000001
$ws.Range("H131").Value = 206009.5
$ws.Range("I131").Value = 334732.34
$ws.Range("K131").Value = 1004197.02
$ws.Range("M131").Value = -999157.02
$ws.Range("H138").Value = 5070.643
$ws.Range("I138").Value = 2332.111
$ws.Range("K138").Value = 6996.333
$ws.Range("M138").Value = -1856.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 730010.4399999999
$ws.Range("I80").Value = 732997.3
$ws.Range("J80").Value = 727023.5600000001
$ws.Range("K80").Value = 732997.3
$ws.Range("L80").Value = 727023.5600000001
$ws.Range("M80").Value = -731999.3
$ws.Range("N80").Value = -729019.5600000001
$ws.Range("H83").Value = 730010.4399999999
$ws.Range("I83").Value = 732997.3
$ws.Range("J83").Value = 727023.5600000001
$ws.Range("K83").Value = 3664986.5
$ws.Range("L83").Value = 3635117.8
$ws.Range("M83").Value = -3659994.5
$ws.Range("N83").Value = -3645101.8
$ws.Range("H96").Value = 30000
$ws.Range("I96").Value = 20000
$ws.Range("K96").Value = 20000
$ws.Range("M96").Value = -17254
$ws.Range("H132").Value = 26201.861
$ws.Range("I132").Value = 2473.1936
$ws.Range("K132").Value = 7419.5808
$ws.Range("M132").Value = -4889.5808

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 391694.12
$ws.Range("I7").Value = 839856.0600000001
$ws.Range("J7").Value = 7555.2856
$ws.Range("K7").Value = 839856.0600000001
$ws.Range("L7").Value = 7555.2856
$ws.Range("M7").Value = -839744.0600000001
$ws.Range("N7").Value = -7779.2856
$ws.Range("H40").Value = 5003590.5
$ws.Range("I40").Value = 6668509.5
$ws.Range("K40").Value = 6668509.5
$ws.Range("M40").Value = -6668373.5
$ws.Range("H42").Value = 17494.125
$ws.Range("J42").Value = 17132.572
$ws.Range("L42").Value = 17132.572
$ws.Range("N42").Value = -18258.572
$ws.Range("H49").Value = 17494.125
$ws.Range("J49").Value = 17132.572
$ws.Range("L49").Value = 17132.572
$ws.Range("N49").Value = -17426.572
$ws.Range("H93").Value = 3389
$ws.Range("I93").Value = 2509.4614
$ws.Range("K93").Value = 2509.4614
$ws.Range("M93").Value = -1261.4614
$ws.Range("H126").Value = 391694.12
$ws.Range("I126").Value = 839856.0600000001
$ws.Range("J126").Value = 7555.2856
$ws.Range("K126").Value = 2519568.18
$ws.Range("L126").Value = 22665.8568
$ws.Range("M126").Value = -2517098.18
$ws.Range("N126").Value = -27605.8568
$ws.Range("H132").Value = 6149.9287
$ws.Range("I132").Value = 3942.8572
$ws.Range("K132").Value = 11828.5716
$ws.Range("M132").Value = -9298.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3656
$ws.Range("I81").Value = 3166.6667
$ws.Range("J81").Value = 4634.6665
$ws.Range("K81").Value = 6333.3334
$ws.Range("L81").Value = 9269.333000000001
$ws.Range("M81").Value = -5272.3334
$ws.Range("N81").Value = -11391.333
$ws.Range("H84").Value = 3656
$ws.Range("I84").Value = 3166.6667
$ws.Range("J84").Value = 4634.6665
$ws.Range("K84").Value = 31666.667
$ws.Range("L84").Value = 46346.665
$ws.Range("M84").Value = -26362.667
$ws.Range("N84").Value = -56954.665
$ws.Range("H101").Value = 46999.5
$ws.Range("J101").Value = 46999.5
$ws.Range("L101").Value = 46999.5
$ws.Range("N101").Value = -53489.5
$ws.Range("H132").Value = 16989.63
$ws.Range("I132").Value = 2755.8447
$ws.Range("K132").Value = 8267.534100000001
$ws.Range("M132").Value = -5737.534100000001
$ws.Range("H136").Value = 304013.6
$ws.Range("I136").Value = 419282.28
$ws.Range("J136").Value = 131110.56
$ws.Range("K136").Value = 1257846.84
$ws.Range("L136").Value = 393331.68
$ws.Range("M136").Value = -1255296.84
$ws.Range("N136").Value = -398431.68
